$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Bmp6"
$ws.Cells.Item(2,3).Value = "Acvr2a"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 31.825501
$ws.Cells.Item(2,8).Value = 95.47650300000001
$ws.Cells.Item(2,9).Value = 0.6063608875535647
$ws.Cells.Item(2,10).Value = 0.6063608875535647
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 14.61878266666667
$ws.Cells.Item(2,14).Value = 43.856348
$ws.Cells.Item(2,15).Value = 0.2662829816142094
$ws.Cells.Item(2,16).Value = 0.2662829816142094
$ws.Cells.Item(2,17).Value = 465.2500823767827
$ws.Cells.Item(2,18).Value = 4187.250741391044
$ws.Cells.Item(2,19).Value = 0.1614635850720016
$ws.Cells.Item(2,20).Value = 0.1614635850720016

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Bmp6"
$ws.Cells.Item(3,3).Value = "Acvr2a"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 31.825501
$ws.Cells.Item(3,8).Value = 95.47650300000001
$ws.Cells.Item(3,9).Value = 0.6063608875535647
$ws.Cells.Item(3,10).Value = 0.6063608875535647
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 27.084169
$ws.Cells.Item(3,14).Value = 81.25250700000001
$ws.Cells.Item(3,15).Value = 0.4933415757187404
$ws.Cells.Item(3,16).Value = 0.4933415757187404
$ws.Cells.Item(3,17).Value = 861.9672475936692
$ws.Cells.Item(3,18).Value = 7757.705228343022
$ws.Cells.Item(3,19).Value = 0.2991430357198895
$ws.Cells.Item(3,20).Value = 0.2991430357198895

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Bmp6"
$ws.Cells.Item(4,3).Value = "Acvr2a"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 31.825501
$ws.Cells.Item(4,8).Value = 95.47650300000001
$ws.Cells.Item(4,9).Value = 0.6063608875535647
$ws.Cells.Item(4,10).Value = 0.6063608875535647
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 13.19647366666667
$ws.Cells.Item(4,14).Value = 39.589421
$ws.Cells.Item(4,15).Value = 0.2403754426670501
$ws.Cells.Item(4,16).Value = 0.2403754426670501
$ws.Cells.Item(4,17).Value = 419.9843858749737
$ws.Cells.Item(4,18).Value = 3779.859472874763
$ws.Cells.Item(4,19).Value = 0.1457542667616735
$ws.Cells.Item(4,20).Value = 0.1457542667616735

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Bmp6"
$ws.Cells.Item(5,3).Value = "Acvr2a"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 4.846280666666666
$ws.Cells.Item(5,8).Value = 14.538842
$ws.Cells.Item(5,9).Value = 0.09233460445363234
$ws.Cells.Item(5,10).Value = 0.09233460445363234
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 14.61878266666667
$ws.Cells.Item(5,14).Value = 43.856348
$ws.Cells.Item(5,15).Value = 0.2662829816142094
$ws.Cells.Item(5,16).Value = 0.2662829816142094
$ws.Cells.Item(5,17).Value = 70.84672380766844
$ws.Cells.Item(5,18).Value = 637.6205142690159
$ws.Cells.Item(5,19).Value = 0.02458713378008188
$ws.Cells.Item(5,20).Value = 0.02458713378008188

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Bmp6"
$ws.Cells.Item(6,3).Value = "Acvr2a"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 4.846280666666666
$ws.Cells.Item(6,8).Value = 14.538842
$ws.Cells.Item(6,9).Value = 0.09233460445363234
$ws.Cells.Item(6,10).Value = 0.09233460445363234
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 27.084169
$ws.Cells.Item(6,14).Value = 81.25250700000001
$ws.Cells.Item(6,15).Value = 0.4933415757187404
$ws.Cells.Item(6,16).Value = 0.4933415757187404
$ws.Cells.Item(6,17).Value = 131.2574845974327
$ws.Cells.Item(6,18).Value = 1181.317361376894
$ws.Cells.Item(6,19).Value = 0.0455524992545216
$ws.Cells.Item(6,20).Value = 0.0455524992545216

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Bmp6"
$ws.Cells.Item(7,3).Value = "Acvr2a"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 4.846280666666666
$ws.Cells.Item(7,8).Value = 14.538842
$ws.Cells.Item(7,9).Value = 0.09233460445363234
$ws.Cells.Item(7,10).Value = 0.09233460445363234
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 13.19647366666667
$ws.Cells.Item(7,14).Value = 39.589421
$ws.Cells.Item(7,15).Value = 0.2403754426670501
$ws.Cells.Item(7,16).Value = 0.2403754426670501
$ws.Cells.Item(7,17).Value = 63.95381519894244
$ws.Cells.Item(7,18).Value = 575.584336790482
$ws.Cells.Item(7,19).Value = 0.02219497141902885
$ws.Cells.Item(7,20).Value = 0.02219497141902885

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Bmp6"
$ws.Cells.Item(8,3).Value = "Acvr2a"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 15.81429
$ws.Cells.Item(8,8).Value = 47.44287
$ws.Cells.Item(8,9).Value = 0.301304507992803
$ws.Cells.Item(8,10).Value = 0.3013045079928031
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 14.61878266666667
$ws.Cells.Item(8,14).Value = 43.856348
$ws.Cells.Item(8,15).Value = 0.2662829816142094
$ws.Cells.Item(8,16).Value = 0.2662829816142094
$ws.Cells.Item(8,17).Value = 231.18566853764
$ws.Cells.Item(8,18).Value = 2080.67101683876
$ws.Cells.Item(8,19).Value = 0.08023226276212599
$ws.Cells.Item(8,20).Value = 0.080232262762126

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Bmp6"
$ws.Cells.Item(9,3).Value = "Acvr2a"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 15.81429
$ws.Cells.Item(9,8).Value = 47.44287
$ws.Cells.Item(9,9).Value = 0.301304507992803
$ws.Cells.Item(9,10).Value = 0.3013045079928031
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 27.084169
$ws.Cells.Item(9,14).Value = 81.25250700000001
$ws.Cells.Item(9,15).Value = 0.4933415757187404
$ws.Cells.Item(9,16).Value = 0.4933415757187404
$ws.Cells.Item(9,17).Value = 428.31690297501
$ws.Cells.Item(9,18).Value = 3854.85212677509
$ws.Cells.Item(9,19).Value = 0.1486460407443292
$ws.Cells.Item(9,20).Value = 0.1486460407443293

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Bmp6"
$ws.Cells.Item(10,3).Value = "Acvr2a"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 15.81429
$ws.Cells.Item(10,8).Value = 47.44287
$ws.Cells.Item(10,9).Value = 0.301304507992803
$ws.Cells.Item(10,10).Value = 0.3013045079928031
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 13.19647366666667
$ws.Cells.Item(10,14).Value = 39.589421
$ws.Cells.Item(10,15).Value = 0.2403754426670501
$ws.Cells.Item(10,16).Value = 0.2403754426670501
$ws.Cells.Item(10,17).Value = 208.69286154203
$ws.Cells.Item(10,18).Value = 1878.23575387827
$ws.Cells.Item(10,19).Value = 0.07242620448634776
$ws.Cells.Item(10,20).Value = 0.07242620448634778
